# Config.xlsx update: introduce MaxContinuousRetryNumber setting and
# refresh the MaxRetryNumber / MaxInitRetryNumber descriptions, then
# leave the workbook focused on the "Constants" sheet (where the new
# setting lives) instead of "Workblocks".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# Insert a new row right after MaxInitRetryNumber (row 4) for the new
# MaxContinuousRetryNumber constant. This pushes every following row
# down by one, including the trailing formatting-only rows.
$ws.Rows.Item(5).Insert()

# Refine the retry-related descriptions to clarify the new local vs.
# Orchestrator queue retry semantics.
$ws.Cells.Item(3, 3).Value = "If > 0, the robot will retry the same transaction which failed with application exception. This is a local data retry. Orchestrator Queue Item retry are managed at the queue level. Must be integer"
$ws.Cells.Item(4, 3).Value = "If > 0 will retry the Initialisation state with a failed exception. Must be an integer."

$ws.Cells.Item(5, 1).Value = "MaxContinuousRetryNumber"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = "If > 0 will keep a record of consecutive failed exceptions of the Process state. When this number is reached, the application will fail. Must be an integer."

# Make "Constants" the active/visible sheet and put the selection on the
# newly inserted row, matching the author's final view state.
$ws.Activate()
$ws.Range("A6").Select()
